$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "20260212_d3c70a24-378f-42e1-b71c-9767071cea25"
$ws.Range("E3").Value = "submitted"
